$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "userUuid": "  ->  "userName": "
# The last edit position (Word's automatic "_GoBack" bookmark) moves
# from the old http://localhost:/8090 split (see Change 2) to right
# after the newly-typed "Name".
# ---------------------------------------------------------------------

$full = $d.Content
$full.Find.Execute("userUuid", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$matchStart = $full.Start
$matchEnd = $full.End

$splitPoint = $matchStart + 4   # position right after "user"

# Temporary bookmark that pins the run boundary between "user" and the
# "Uuid"/"Name" text while we edit, so Word doesn't re-flow/merge the
# "user" run together with the replacement text.
$d.Bookmarks.Add("ZZTemp1", $d.Range($splitPoint, $splitPoint))

# Re-seat "_GoBack" at the boundary right after the new "Name" text.
# Because a bookmark name is unique, this single Add call both creates
# it here *and* removes it from its old location (the localhost/8090
# split handled in Change 2 below).
$d.Bookmarks.Add("_GoBack", $d.Range($matchEnd, $matchEnd))

$colonEnd = $matchEnd + 4
# Temporary bookmark pinning the boundary right after the closing
# quote/colon text so it doesn't re-merge with the following value.
$d.Bookmarks.Add("ZZTemp2", $d.Range($colonEnd, $colonEnd))

# Replace "Uuid" with "Name"
$r = $d.Range($splitPoint, $matchEnd)
$r.Text = "Name"

# Re-touch the `": "` piece so it drops its stale formatting id and
# becomes its own clean run bordered by the bookmark on one side and
# the untouched value run on the other.
$colonRange = $d.Range($matchEnd, $colonEnd)
$colonRange.Text = "zzzz"
$colonRange2 = $d.Range($matchEnd, $matchEnd + 4)
$colonRange2.Text = """: """

# Drop the temporary boundary markers, leaving "_GoBack" in its new spot.
$d.Bookmarks("ZZTemp1").Delete()
$d.Bookmarks("ZZTemp2").Delete()

# ---------------------------------------------------------------------
# Change 2: http://localhost: | 8090/crimepatrol/.../aggregatedinfo
# merge back into a single run (the bookmark that used to split them
# was relocated above).
# ---------------------------------------------------------------------

$urlFind = $d.Content
$urlFind.Find.Execute("aggregatedinfo", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$urlEnd = $urlFind.End
$urlStart = $urlEnd - 72   # length of "http://localhost:8090/crimepatrol/expenditures/Settlement/aggregatedinfo"

# Keep the "[GET] " run untouched by pinning the boundary in front of
# the URL while we rebuild it as a single run.
$d.Bookmarks.Add("ZZTemp3", $d.Range($urlStart, $urlStart))

$urlRange = $d.Range($urlStart, $urlEnd)
$urlRange.Text = "zzzz"
$urlRange2 = $d.Range($urlStart, $urlStart + 4)
$urlRange2.Text = "http://localhost:8090/crimepatrol/expenditures/Settlement/aggregatedinfo"

$d.Bookmarks("ZZTemp3").Delete()
